$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in row 1: I1 = "I0", J1 = "IF"
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the existing header formatting (bold, bordered, centered) used by
# B1:H1 -- copy it from H1 onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New data cells in row 2: I2 = 8, J2 = 8
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
